# Update for popup added
# Append rows 3..9 of nurse data below the existing header (row 1) and
# first data row (row 2), matching the "Daniel" row's formatting:
#   - A, E : plain numbers (no special formatting)
#   - B    : plain text (name)
#   - C    : text-typed numeric string (license number stored as text,
#            same as the existing C2 cell) with NO explicit cell style
#   - D    : date/time serial number carrying the same cell style as D2
#            (numFmtId 14 date format, style index reused from D2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats / xlPasteValues constants used below
$xlPasteFormats = -4122
$xlPasteValues  = -4163

function Add-NurseRow {
    param([int]$Row, [int]$Id, [string]$Name, [string]$License, [double]$Dob, [int]$Age)

    # A: id (plain number)
    $ws.Range("A$Row").Value = $Id

    # B: name (plain text)
    $ws.Range("B$Row").Value = $Name

    # C: license number stored as TEXT (mirrors source row 2, which has
    # t="str" / no explicit style). Assigning the digit string directly
    # would auto-convert to a number, so we route it through a
    # string-literal formula and then freeze it back to a plain value via
    # Copy / PasteSpecial(xlPasteValues) -- this collapses the formula to
    # a literal string cell without adding a NumberFormat-driven style.
    $ws.Range("C$Row").Formula = '="' + $License + '"'
    $ws.Range("C$Row").Copy() | Out-Null
    $ws.Range("C$Row").PasteSpecial($xlPasteValues) | Out-Null

    # D: date of birth serial, re-using D2's exact style (numFmtId 14)
    # by copying its format only, then writing the numeric value.
    $ws.Range("D2").Copy() | Out-Null
    $ws.Range("D$Row").PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("D$Row").Value = $Dob

    # E: age (plain number)
    $ws.Range("E$Row").Value = $Age
}

Add-NurseRow 3 3  "Daniel"  "2147483647"  39008.00011574074    18
Add-NurseRow 4 4  "Passion" "1231234324"  39009.00011574074    18
Add-NurseRow 5 5  "Kumar"   "2147483647"  39008.00011574074    18
Add-NurseRow 6 6  "Kamal"   "2147483647"  39008.00011574074    18
Add-NurseRow 7 7  "Suresh"  "2147483647"  32434.000115740742   35
Add-NurseRow 8 9  "qwwerr"  "7854551266"  37356.00011574074    23
Add-NurseRow 9 10 "ashish"  "78955952599" 36642.00011574074    25

# Clear the leftover clipboard/marching-ants reference left by Copy().
$excel.CutCopyMode = 0
